# Auto-generated edit script for Kujata_Profits workbook update
# Updates market-price / profit columns (H..N) on several rows across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 6122.9165
$ws.Cells.Item(98, 9).Value = 7831.1113
$ws.Cells.Item(98, 10).Value = 998.3333
$ws.Cells.Item(98, 11).Value = 7831.1113
$ws.Cells.Item(98, 12).Value = 998.3333
$ws.Cells.Item(98, 13).Value = -6333.1113
$ws.Cells.Item(98, 14).Value = -3994.3333

$ws.Cells.Item(122, 8).Value = 6122.9165
$ws.Cells.Item(122, 9).Value = 7831.1113
$ws.Cells.Item(122, 10).Value = 998.3333
$ws.Cells.Item(122, 11).Value = 23493.3339
$ws.Cells.Item(122, 12).Value = 2994.9999
$ws.Cells.Item(122, 13).Value = -21043.3339
$ws.Cells.Item(122, 14).Value = -7894.9999

$ws.Cells.Item(132, 8).Value = 16679448
$ws.Cells.Item(132, 9).Value = 33349524
$ws.Cells.Item(132, 10).Value = 9371.1
$ws.Cells.Item(132, 11).Value = 100048572
$ws.Cells.Item(132, 12).Value = 28113.3
$ws.Cells.Item(132, 13).Value = -100046042
$ws.Cells.Item(132, 14).Value = -33173.3

$ws.Cells.Item(135, 8).Value = 1554.65
$ws.Cells.Item(135, 9).Value = 598.5625
$ws.Cells.Item(135, 11).Value = 5387.0625
$ws.Cells.Item(135, 13).Value = -2852.0625

$ws.Cells.Item(137, 8).Value = 1335.1628
$ws.Cells.Item(137, 9).Value = 1028
$ws.Cells.Item(137, 10).Value = 1853.5
$ws.Cells.Item(137, 11).Value = 3084
$ws.Cells.Item(137, 12).Value = 5560.5
$ws.Cells.Item(137, 13).Value = -534
$ws.Cells.Item(137, 14).Value = -10660.5

$ws.Cells.Item(138, 8).Value = 501536.22
$ws.Cells.Item(138, 9).Value = 1998.5
$ws.Cells.Item(138, 10).Value = 581462.25
$ws.Cells.Item(138, 11).Value = 5995.5
$ws.Cells.Item(138, 12).Value = 1744386.75
$ws.Cells.Item(138, 13).Value = -855.5
$ws.Cells.Item(138, 14).Value = -1754666.75

$ws.Cells.Item(141, 8).Value = 3448.2222
$ws.Cells.Item(141, 9).Value = 3629.25
$ws.Cells.Item(141, 11).Value = 10887.75
$ws.Cells.Item(141, 13).Value = -5707.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2169.79
$ws.Cells.Item(32, 9).Value = 1739.1235
$ws.Cells.Item(32, 11).Value = 1739.1235
$ws.Cells.Item(32, 13).Value = -1452.1235

$ws.Cells.Item(61, 8).Value = 76924180
$ws.Cells.Item(61, 9).Value = 83334280
$ws.Cells.Item(61, 11).Value = 83334280
$ws.Cells.Item(61, 13).Value = -83334068

$ws.Cells.Item(74, 8).Value = 2008
$ws.Cells.Item(74, 9).Value = 1522
$ws.Cells.Item(74, 11).Value = 1522
$ws.Cells.Item(74, 13).Value = -648

$ws.Cells.Item(77, 8).Value = 2008
$ws.Cells.Item(77, 9).Value = 1522
$ws.Cells.Item(77, 11).Value = 7610
$ws.Cells.Item(77, 13).Value = -3242

$ws.Cells.Item(132, 8).Value = 1790.9833
$ws.Cells.Item(132, 9).Value = 1384.122
$ws.Cells.Item(132, 10).Value = 2668.9473
$ws.Cells.Item(132, 11).Value = 4152.366
$ws.Cells.Item(132, 12).Value = 8006.841899999999
$ws.Cells.Item(132, 13).Value = -1622.366
$ws.Cells.Item(132, 14).Value = -13066.8419

$ws.Cells.Item(136, 8).Value = 76924180
$ws.Cells.Item(136, 9).Value = 83334280
$ws.Cells.Item(136, 11).Value = 250002840
$ws.Cells.Item(136, 13).Value = -250000290

$ws.Cells.Item(138, 8).Value = 47322.375
$ws.Cells.Item(138, 10).Value = 47322.375
$ws.Cells.Item(138, 12).Value = 47322.375
$ws.Cells.Item(138, 14).Value = -57602.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 415.18182
$ws.Cells.Item(22, 9).Value = 190
$ws.Cells.Item(22, 10).Value = 465.22223
$ws.Cells.Item(22, 11).Value = 190
$ws.Cells.Item(22, 12).Value = 465.22223
$ws.Cells.Item(22, 13).Value = -17
$ws.Cells.Item(22, 14).Value = -811.2222300000001

$ws.Cells.Item(64, 8).Value = 500.58334
$ws.Cells.Item(64, 9).Value = 343.85715
$ws.Cells.Item(64, 10).Value = 720
$ws.Cells.Item(64, 11).Value = 343.85715
$ws.Cells.Item(64, 12).Value = 720
$ws.Cells.Item(64, 13).Value = -118.85715
$ws.Cells.Item(64, 14).Value = -1170

$ws.Cells.Item(67, 8).Value = 500.58334
$ws.Cells.Item(67, 9).Value = 343.85715
$ws.Cells.Item(67, 10).Value = 720
$ws.Cells.Item(67, 11).Value = 343.85715
$ws.Cells.Item(67, 12).Value = 720
$ws.Cells.Item(67, 13).Value = 436.14285
$ws.Cells.Item(67, 14).Value = -2280

$ws.Cells.Item(139, 8).Value = 39290
$ws.Cells.Item(139, 10).Value = 39290
$ws.Cells.Item(139, 12).Value = 39290
$ws.Cells.Item(139, 14).Value = -49570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1781.1515
$ws.Cells.Item(31, 9).Value = 1619.5555
$ws.Cells.Item(31, 11).Value = 1619.5555
$ws.Cells.Item(31, 13).Value = -1324.5555

$ws.Cells.Item(34, 8).Value = 1781.1515
$ws.Cells.Item(34, 9).Value = 1619.5555
$ws.Cells.Item(34, 11).Value = 1619.5555
$ws.Cells.Item(34, 13).Value = -1417.5555

$ws.Cells.Item(58, 8).Value = 1591.6945
$ws.Cells.Item(58, 9).Value = 1217.2727
$ws.Cells.Item(58, 11).Value = 1217.2727
$ws.Cells.Item(58, 13).Value = -1014.2727

$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Range("N133").ClearContents()

$ws.Cells.Item(134, 8).Value = 50002300
$ws.Cells.Item(134, 9).Value = 2583.1428
$ws.Cells.Item(134, 10).Value = 166668300
$ws.Cells.Item(134, 11).Value = 7749.428400000001
$ws.Cells.Item(134, 12).Value = 500004900
$ws.Cells.Item(134, 13).Value = -5214.428400000001
$ws.Cells.Item(134, 14).Value = -500009970

$ws.Cells.Item(136, 8).Value = 1591.6945
$ws.Cells.Item(136, 9).Value = 1217.2727
$ws.Cells.Item(136, 11).Value = 3651.8181
$ws.Cells.Item(136, 13).Value = -1101.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 266
$ws.Cells.Item(11, 9).Value = 266
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 798
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -658
$ws.Range("N11").ClearContents()

$ws.Cells.Item(100, 8).Value = 15106.875
$ws.Cells.Item(100, 10).Value = 15106.875
$ws.Cells.Item(100, 12).Value = 45320.625
$ws.Cells.Item(100, 14).Value = -46942.625

$ws.Cells.Item(131, 8).Value = 16394353
$ws.Cells.Item(131, 9).Value = 55556044
$ws.Cells.Item(131, 10).Value = 1087.2094
$ws.Cells.Item(131, 11).Value = 166668132
$ws.Cells.Item(131, 12).Value = 3261.6282
$ws.Cells.Item(131, 13).Value = -166663092
$ws.Cells.Item(131, 14).Value = -13341.6282

$ws.Cells.Item(141, 8).Value = 90912140
$ws.Cells.Item(141, 9).Value = 111112904
$ws.Cells.Item(141, 11).Value = 333338712
$ws.Cells.Item(141, 13).Value = -333333532

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 5100
$ws.Cells.Item(18, 10).Value = 5100
$ws.Cells.Item(18, 12).Value = 5100
$ws.Cells.Item(18, 14).Value = -5686

$ws.Cells.Item(80, 8).Value = 4690.6924
$ws.Cells.Item(80, 9).Value = 4072.375
$ws.Cells.Item(80, 10).Value = 5680
$ws.Cells.Item(80, 11).Value = 4072.375
$ws.Cells.Item(80, 12).Value = 5680
$ws.Cells.Item(80, 13).Value = -3074.375
$ws.Cells.Item(80, 14).Value = -7676

$ws.Cells.Item(83, 8).Value = 4690.6924
$ws.Cells.Item(83, 9).Value = 4072.375
$ws.Cells.Item(83, 10).Value = 5680
$ws.Cells.Item(83, 11).Value = 20361.875
$ws.Cells.Item(83, 12).Value = 28400
$ws.Cells.Item(83, 13).Value = -15369.875
$ws.Cells.Item(83, 14).Value = -38384

$ws.Cells.Item(126, 8).Value = 1745.4667
$ws.Cells.Item(126, 9).Value = 1533.5
$ws.Cells.Item(126, 11).Value = 4600.5
$ws.Cells.Item(126, 13).Value = -2130.5

$ws.Cells.Item(132, 8).Value = 2679.3142
$ws.Cells.Item(132, 9).Value = 2446.5217
$ws.Cells.Item(132, 11).Value = 7339.5651
$ws.Cells.Item(132, 13).Value = -4809.5651

$ws.Cells.Item(139, 8).Value = 31413.2
$ws.Cells.Item(139, 10).Value = 31413.2
$ws.Cells.Item(139, 12).Value = 31413.2
$ws.Cells.Item(139, 14).Value = -41693.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 1088.5
$ws.Cells.Item(9, 9).Value = 451.33334
$ws.Cells.Item(9, 11).Value = 451.33334
$ws.Cells.Item(9, 13).Value = -227.33334

$ws.Cells.Item(16, 8).Value = 895.3570999999999
$ws.Cells.Item(16, 9).Value = 895.3570999999999
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 895.3570999999999
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -725.3570999999999
$ws.Range("N16").ClearContents()

$ws.Cells.Item(82, 8).Value = 1630.75
$ws.Cells.Item(82, 9).Value = 1576
$ws.Cells.Item(82, 10).Value = 1795
$ws.Cells.Item(82, 11).Value = 1576
$ws.Cells.Item(82, 12).Value = 1795
$ws.Cells.Item(82, 13).Value = -1215
$ws.Cells.Item(82, 14).Value = -2517

$ws.Cells.Item(85, 8).Value = 1630.75
$ws.Cells.Item(85, 9).Value = 1576
$ws.Cells.Item(85, 10).Value = 1795
$ws.Cells.Item(85, 11).Value = 1576
$ws.Cells.Item(85, 12).Value = 1795
$ws.Cells.Item(85, 13).Value = -328
$ws.Cells.Item(85, 14).Value = -4291

$ws.Cells.Item(136, 8).Value = 1162.7307
$ws.Cells.Item(136, 9).Value = 1009.7083
$ws.Cells.Item(136, 11).Value = 3029.1249
$ws.Cells.Item(136, 13).Value = -479.1248999999998

$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 55558356
$ws.Cells.Item(62, 10).Value = 3239.6
$ws.Cells.Item(62, 12).Value = 3239.6
$ws.Cells.Item(62, 14).Value = -4487.6

$ws.Cells.Item(65, 8).Value = 55558356
$ws.Cells.Item(65, 10).Value = 3239.6
$ws.Cells.Item(65, 12).Value = 16198
$ws.Cells.Item(65, 14).Value = -22438

$ws.Cells.Item(136, 8).Value = 1746.8572
$ws.Cells.Item(136, 9).Value = 1520.7693
$ws.Cells.Item(136, 11).Value = 4562.3079
$ws.Cells.Item(136, 13).Value = -2012.3079
